$wb = $excel.ActiveWorkbook

# ===== Sheet 1 =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 09:21:49"
$ws.Range("A3").Value = "Total filas: 125"

$rows1 = @(
    @("05:42:22", "05:43", "14_ABASTO", 1, "LP1912"),
    @("05:42:22", "05:52", "17_ROMERO", 10, "LP1912"),
    @("05:42:22", "06:01", "16_SANTA ANA", 19, "LP1912"),
    @("05:42:22", "06:04", "10_OLMOS", 22, "LP1912"),
    @("05:42:22", "06:11", "215A_EL PATO", 29, "LP1912"),
    @("05:42:22", "06:24", "11_ETCHEVERRY", 42, "LP1912"),
    @("05:42:22", "06:27", "23_HERNANDEZ", 45, "LP1912"),
    @("05:42:22", "06:31", "17X38_ROMERO", 49, "LP1912"),
    @("05:42:22", "06:31", "16_SANTA ANA", 49, "LP1912"),
    @("05:42:22", "06:39", "225_C ROCA-H SUR", 57, "LP1912"),
    @("06:33:46", "06:50", "215A_EL PATO", 17, "LP1912"),
    @("06:16:15", "06:51", "215A_EL PATO", 35, "LP1912"),
    @("06:52:52", "06:52", "215A_EL PATO", 0, "LP1912"),
    @("05:42:22", "06:54", "14_ABASTO", 72, "LP1912"),
    @("06:52:52", "06:55", "14_ABASTO", 3, "LP1912"),
    @("06:45:50", "07:00", "16_SANTA ANA", 15, "LP1912"),
    @("05:57:08", "07:01", "16_SANTA ANA", 64, "LP1912"),
    @("06:45:50", "07:03", "225_GOMEZ", 18, "LP1912"),
    @("05:42:22", "07:04", "225_GOMEZ", 82, "LP1912"),
    @("06:33:46", "07:06", "215C_EL PATO", 33, "LP1912"),
    @("05:42:22", "07:07", "215C_EL PATO", 85, "LP1912"),
    @("06:33:46", "07:13", "14X44_ABASTO", 40, "LP1912"),
    @("05:42:22", "07:14", "14X44_ABASTO", 92, "LP1912"),
    @("06:33:46", "07:15", "16_SANTA ANA", 42, "LP1912"),
    @("06:33:46", "07:20", "215A_EL PATO", 47, "LP1912"),
    @("05:42:22", "07:21", "215A_EL PATO", 99, "LP1912"),
    @("06:52:52", "07:28", "14_ABASTO", 36, "LP1912"),
    @("05:57:08", "07:29", "14_ABASTO", 92, "LP1912"),
    @("05:42:22", "07:33", "23_HERNANDEZ", 111, "LP1912"),
    @("05:57:08", "07:34", "23_HERNANDEZ", 97, "LP1912"),
    @("06:45:50", "07:35", "17X38_ROMERO", 50, "LP1912"),
    @("05:42:22", "07:36", "17X38_ROMERO", 114, "LP1912"),
    @("06:33:46", "07:36", "27_EL RETIRO", 63, "LP1912"),
    @("05:42:22", "07:37", "27_EL RETIRO", 115, "LP1912"),
    @("07:36:59", "07:38", "16_SANTA ANA", 2, "LP1912"),
    @("06:33:46", "07:43", "10_OLMOS", 70, "LP1912"),
    @("05:57:08", "07:44", "10_OLMOS", 107, "LP1912"),
    @("07:36:59", "07:50", "15_ABASTO", 14, "LP1912"),
    @("05:57:08", "07:51", "15_ABASTO", 114, "LP1912"),
    @("06:16:15", "07:58", "23_HERNANDEZ", 102, "LP1912"),
    @("07:12:53", "07:59", "23_HERNANDEZ", 47, "LP1912"),
    @("06:33:46", "07:59", "11_ETCHEVERRY", 86, "LP1912"),
    @("06:33:46", "08:00", "23_HERNANDEZ", 87, "LP1912"),
    @("06:16:15", "08:00", "11_ETCHEVERRY", 104, "LP1912"),
    @("06:45:50", "08:01", "23_HERNANDEZ", 76, "LP1912"),
    @("07:12:53", "08:01", "16_SANTA ANA", 49, "LP1912"),
    @("06:52:52", "08:02", "23_HERNANDEZ", 70, "LP1912"),
    @("07:36:59", "08:02", "17_ROMERO", 26, "LP1912"),
    @("06:16:15", "08:03", "17_ROMERO", 107, "LP1912"),
    @("06:33:46", "08:12", "10_OLMOS", 99, "LP1912"),
    @("07:12:53", "08:13", "10_OLMOS", 61, "LP1912"),
    @("08:11:27", "08:14", "11_ETCHEVERRY", 3, "LP1912"),
    @("07:48:35", "08:14", "10_OLMOS", 26, "LP1912"),
    @("07:36:59", "08:14", "17_ROMERO", 38, "LP1912"),
    @("06:16:15", "08:15", "17_ROMERO", 119, "LP1912"),
    @("07:36:59", "08:25", "15X38_ABASTO", 49, "LP1912"),
    @("06:33:46", "08:26", "15X38_ABASTO", 113, "LP1912"),
    @("06:33:46", "08:27", "84_COLONIA URQUIZA-ESC 49", 114, "LP1912"),
    @("08:29:19", "08:29", "23_HERNANDEZ", 0, "LP1912"),
    @("06:45:50", "08:29", "14_ABASTO", 104, "LP1912"),
    @("07:36:59", "08:30", "16_P MOR-SANTA ANA", 54, "LP1912"),
    @("06:33:46", "08:31", "16_P MOR-SANTA ANA", 118, "LP1912"),
    @("08:11:27", "08:33", "23_HERNANDEZ", 22, "LP1912"),
    @("06:45:50", "08:38", "215C_EL PATO", 113, "LP1912"),
    @("07:48:35", "08:39", "215C_EL PATO", 51, "LP1912"),
    @("07:12:53", "08:43", "10_OLMOS", 91, "LP1912"),
    @("07:48:35", "08:44", "10_OLMOS", 56, "LP1912"),
    @("07:12:53", "08:49", "215A_EL PATO", 97, "LP1912"),
    @("07:48:35", "08:50", "215A_EL PATO", 62, "LP1912"),
    @("08:11:27", "08:53", "16_SANTA ANA", 42, "LP1912"),
    @("08:29:19", "08:54", "16_SANTA ANA", 25, "LP1912"),
    @("08:53:12", "08:55", "16_SANTA ANA", 2, "LP1912"),
    @("07:12:53", "08:59", "215B_EL PATO", 107, "LP1912"),
    @("08:39:08", "09:00", "215B_EL PATO", 21, "LP1912"),
    @("07:36:59", "09:01", "17X38_ROMERO", 85, "LP1912"),
    @("07:36:59", "09:02", "23_HERNANDEZ", 86, "LP1912"),
    @("07:12:53", "09:02", "17X38_ROMERO", 110, "LP1912"),
    @("08:29:19", "09:03", "23_HERNANDEZ", 34, "LP1912"),
    @("08:39:08", "09:04", "23_HERNANDEZ", 25, "LP1912"),
    @("07:36:59", "09:04", "16_SANTA ANA", 88, "LP1912"),
    @("07:48:35", "09:08", "16_SANTA ANA", 80, "LP1912"),
    @("07:36:59", "09:14", "15_ABASTO", 98, "LP1912"),
    @("07:36:59", "09:14", "11_ETCHEVERRY", 98, "LP1912"),
    @("08:39:08", "09:15", "11_ETCHEVERRY", 36, "LP1912"),
    @("07:36:59", "09:16", "27_EL RETIRO", 100, "LP1912"),
    @("08:39:08", "09:17", "27_EL RETIRO", 38, "LP1912"),
    @("07:55:46", "09:21", "16_SANTA ANA", 86, "LP1912"),
    @("09:21:49", "09:22", "16_SANTA ANA", 1, "LP1912"),
    @("07:36:59", "09:26", "215_EL PELIGRO", 110, "LP1912"),
    @("07:48:35", "09:27", "215_EL PELIGRO", 99, "LP1912"),
    @("07:36:59", "09:30", "16_P MOR-SANTA ANA", 114, "LP1912"),
    @("07:48:35", "09:31", "16_P MOR-SANTA ANA", 103, "LP1912"),
    @("08:39:08", "09:32", "23_HERNANDEZ", 53, "LP1912"),
    @("08:46:25", "09:33", "23_HERNANDEZ", 47, "LP1912"),
    @("07:48:35", "09:39", "15_ABASTO", 111, "LP1912"),
    @("07:48:35", "09:44", "14_ABASTO", 116, "LP1912"),
    @("08:39:08", "09:45", "14_ABASTO", 66, "LP1912"),
    @("07:55:46", "09:51", "15_ABASTO", 116, "LP1912"),
    @("08:53:12", "09:54", "10_OLMOS", 61, "LP1912"),
    @("09:21:49", "10:01", "16_SANTA ANA", 40, "LP1912"),
    @("09:21:49", "10:02", "215C_EL PATO", 41, "LP1912"),
    @("08:11:27", "10:03", "215C_EL PATO", 112, "LP1912"),
    @("09:21:49", "10:03", "23_HERNANDEZ", 42, "LP1912"),
    @("08:46:25", "10:04", "14_ABASTO", 78, "LP1912"),
    @("08:39:08", "10:05", "14_ABASTO", 86, "LP1912"),
    @("08:11:27", "10:10", "10_OLMOS", 119, "LP1912"),
    @("08:29:19", "10:11", "27_EL RETIRO", 102, "LP1912"),
    @("08:29:19", "10:12", "10_OLMOS", 103, "LP1912"),
    @("09:21:49", "10:14", "17_ROMERO", 53, "LP1912"),
    @("08:29:19", "10:14", "10_OLMOS", 105, "LP1912"),
    @("08:29:19", "10:15", "17_ROMERO", 106, "LP1912"),
    @("08:39:08", "10:15", "10_OLMOS", 96, "LP1912"),
    @("08:46:25", "10:16", "10_OLMOS", 90, "LP1912"),
    @("08:46:25", "10:18", "27_EL RETIRO", 92, "LP1912"),
    @("09:21:49", "10:24", "11_ETCHEVERRY", 63, "LP1912"),
    @("08:29:19", "10:26", "15X38_ABASTO", 117, "LP1912"),
    @("08:39:08", "10:30", "11_ETCHEVERRY", 111, "LP1912"),
    @("08:39:08", "10:34", "10_OLMOS", 115, "LP1912"),
    @("08:39:08", "10:37", "16_P MOR-SANTA ANA", 118, "LP1912"),
    @("09:21:49", "10:44", "14_ABASTO", 83, "LP1912"),
    @("09:21:49", "10:56", "27_EL RETIRO", 95, "LP1912"),
    @("09:21:49", "11:01", "17_ROMERO", 100, "LP1912"),
    @("09:21:49", "11:04", "14_ABASTO", 103, "LP1912"),
    @("09:21:49", "11:14", "225_C ROCA-H SUR", 113, "LP1912"),
    @("09:21:49", "11:20", "215C_EL PATO", 119, "LP1912")
)
$r = 6
foreach ($row in $rows1) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ===== Sheet 2 =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 09:21:49"
$ws.Range("A3").Value = "Total filas: 19"

$rows2 = @(
    @("05:42:22", "06:11", "215A_EL PATO", 29, "LP1912"),
    @("06:33:46", "06:50", "215A_EL PATO", 17, "LP1912"),
    @("06:16:15", "06:51", "215A_EL PATO", 35, "LP1912"),
    @("06:52:52", "06:52", "215A_EL PATO", 0, "LP1912"),
    @("06:33:46", "07:06", "215C_EL PATO", 33, "LP1912"),
    @("05:42:22", "07:07", "215C_EL PATO", 85, "LP1912"),
    @("06:33:46", "07:20", "215A_EL PATO", 47, "LP1912"),
    @("05:42:22", "07:21", "215A_EL PATO", 99, "LP1912"),
    @("06:45:50", "08:38", "215C_EL PATO", 113, "LP1912"),
    @("07:48:35", "08:39", "215C_EL PATO", 51, "LP1912"),
    @("07:12:53", "08:49", "215A_EL PATO", 97, "LP1912"),
    @("07:48:35", "08:50", "215A_EL PATO", 62, "LP1912"),
    @("07:12:53", "08:59", "215B_EL PATO", 107, "LP1912"),
    @("08:39:08", "09:00", "215B_EL PATO", 21, "LP1912"),
    @("07:36:59", "09:26", "215_EL PELIGRO", 110, "LP1912"),
    @("07:48:35", "09:27", "215_EL PELIGRO", 99, "LP1912"),
    @("09:21:49", "10:02", "215C_EL PATO", 41, "LP1912"),
    @("08:11:27", "10:03", "215C_EL PATO", 112, "LP1912"),
    @("09:21:49", "11:20", "215C_EL PATO", 119, "LP1912")
)
$r = 6
foreach ($row in $rows2) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ===== Sheet 3 =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 09:21:49"
$ws.Range("A3").Value = "Total filas: 18"

$rows3 = @(
    @("05:42:22", "07:27", "215A_LA PLATA", 105, "L6173"),
    @("06:33:46", "08:09", "215A_LA PLATA", 96, "L6173"),
    @("06:16:15", "08:10", "215A_LA PLATA", 114, "L6173"),
    @("08:11:27", "08:11", "215A_LA PLATA", 0, "L6173"),
    @("06:33:46", "08:22", "215C_LA PLATA", 109, "L6203"),
    @("07:48:35", "08:25", "215C_LA PLATA", 37, "L6203"),
    @("07:55:46", "08:26", "215C_LA PLATA", 31, "L6203"),
    @("07:36:59", "08:27", "215C_LA PLATA", 51, "L6203"),
    @("08:46:25", "08:48", "215A_LA PLATA", 2, "L6173"),
    @("07:36:59", "08:51", "215A_LA PLATA", 75, "L6173"),
    @("07:48:35", "08:52", "215A_LA PLATA", 64, "L6173"),
    @("08:53:12", "08:53", "215A_LA PLATA", 0, "L6173"),
    @("09:21:49", "10:08", "215C_LA PLATA", 47, "L6203"),
    @("08:11:27", "10:09", "215C_LA PLATA", 118, "L6203"),
    @("09:21:49", "10:22", "215A_LA PLATA", 61, "L6173"),
    @("08:29:19", "10:23", "215A_LA PLATA", 114, "L6173"),
    @("09:21:49", "10:30", "215B_LP-P MOR-1 Y 57", 69, "L6173"),
    @("08:39:08", "10:31", "215B_LP-P MOR-1 Y 57", 112, "L6173")
)
$r = 6
foreach ($row in $rows3) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
